# Weekly fruit/vegetable price refresh for "Hortaliza, Terminal
# Hortofrutícola Agro Chillán - Cebollín" (commit: "Fruta / hortaliza,
# semanal"). The price-history rows 196-230 are refreshed with the next
# reporting cycle's values, and two brand-new weekly rows (231-232) are
# appended, growing the sheet from A1:R230 to A1:R232.
#
# Each inner array is: RowNumber, A..R (18 columns) in the same column
# order as the sheet header (Mercado ID, Mercado, Región, Fecha, Codreg,
# Categoría ID, Categoría, Variedad, Calidad, Volumen, Precio mínimo,
# Precio máximo, Precio promedio ponderado, Unidad de comercialización,
# Origen, Precio $/Kg, Kg o Unidades, Clasificación). Column D (Fecha) is
# stored as the underlying Excel date serial number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @(196,7,"Terminal Hortofrutícola Agro Chillán","Ñuble",45244,16,100112037,"Cebollín","Sin especificar","Primera",300,6000,6000,6000,"$/paquete 36 unidades","Región Metropolitana",167,36,"Hortaliza"),
    @(197,7,"Terminal Hortofrutícola Agro Chillán","Ñuble",45244,16,100112037,"Cebollín","Sin especificar","Segunda",300,5000,5000,5000,"$/paquete 36 unidades","Región Metropolitana",139,36,"Hortaliza"),
    @(198,7,"Terminal Hortofrutícola Agro Chillán","Ñuble",45195,16,100112037,"Cebollín","Sin especificar","Primera",200,6500,6500,6500,"$/paquete 36 unidades","Provincia de Diguillín",181,36,"Hortaliza"),
    @(199,7,"Terminal Hortofrutícola Agro Chillán","Ñuble",45090,16,100112037,"Cebollín","Sin especificar","Primera",150,6000,6000,6000,"$/paquete 36 unidades","Provincia de Diguillín",167,36,"Hortaliza"),
    @(200,7,"Terminal Hortofrutícola Agro Chillán","Ñuble",45090,16,100112037,"Cebollín","Sin especificar","Segunda",100,5000,5000,5000,"$/paquete 36 unidades","Provincia de Diguillín",139,36,"Hortaliza"),
    @(201,7,"Terminal Hortofrutícola Agro Chillán","Ñuble",44847,16,100112037,"Cebollín","Sin especificar","Primera",120,7500,8000,7750,"$/docena de atados","Provincia de Diguillín",2583,3,"Hortaliza"),
    @(202,7,"Terminal Hortofrutícola Agro Chillán","Ñuble",45209,16,100112037,"Cebollín","Sin especificar","Primera",150,6000,6000,6000,"$/paquete 36 unidades","Provincia de Diguillín",167,36,"Hortaliza"),
    @(203,7,"Terminal Hortofrutícola Agro Chillán","Ñuble",45040,16,100112037,"Cebollín","Sin especificar","Primera",60,7000,7000,7000,"$/paquete 36 unidades","Provincia de Diguillín",194,36,"Hortaliza"),
    @(204,7,"Terminal Hortofrutícola Agro Chillán","Ñuble",45168,16,100112037,"Cebollín","Sin especificar","Primera",120,6000,6000,6000,"$/paquete 36 unidades","Provincia de Diguillín",167,36,"Hortaliza"),
    @(205,7,"Terminal Hortofrutícola Agro Chillán","Ñuble",45070,16,100112037,"Cebollín","Sin especificar","Primera",80,7000,7000,7000,"$/paquete 36 unidades","Provincia de Diguillín",194,36,"Hortaliza"),
    @(206,7,"Terminal Hortofrutícola Agro Chillán","Ñuble",44985,16,100112037,"Cebollín","Sin especificar","Primera",80,6000,6000,6000,"$/paquete 36 unidades","Provincia de Diguillín",167,36,"Hortaliza"),
    @(207,7,"Terminal Hortofrutícola Agro Chillán","Ñuble",44985,16,100112037,"Cebollín","Sin especificar","Primera",150,800,800,800,"$/paquete 6 unidades","Provincia de Diguillín",133,6,"Hortaliza"),
    @(208,7,"Terminal Hortofrutícola Agro Chillán","Ñuble",44985,16,100112037,"Cebollín","Sin especificar","Segunda",150,600,600,600,"$/paquete 6 unidades","Provincia de Diguillín",100,6,"Hortaliza"),
    @(209,7,"Terminal Hortofrutícola Agro Chillán","Ñuble",44973,16,100112037,"Cebollín","Sin especificar","Primera",300,800,800,800,"$/paquete 6 unidades","Provincia de Diguillín",133,6,"Hortaliza"),
    @(210,7,"Terminal Hortofrutícola Agro Chillán","Ñuble",44973,16,100112037,"Cebollín","Sin especificar","Segunda",200,600,600,600,"$/paquete 6 unidades","Provincia de Diguillín",100,6,"Hortaliza"),
    @(211,7,"Terminal Hortofrutícola Agro Chillán","Ñuble",45219,16,100112037,"Cebollín","Sin especificar","Primera",250,6000,6000,6000,"$/paquete 36 unidades","Provincia de Diguillín",167,36,"Hortaliza"),
    @(212,7,"Terminal Hortofrutícola Agro Chillán","Ñuble",44935,16,100112037,"Cebollín","Sin especificar","Primera",200,600,600,600,"$/paquete 6 unidades","Provincia de Diguillín",100,6,"Hortaliza"),
    @(213,7,"Terminal Hortofrutícola Agro Chillán","Ñuble",44935,16,100112037,"Cebollín","Sin especificar","Segunda",150,500,500,500,"$/paquete 6 unidades","Provincia de Diguillín",83,6,"Hortaliza"),
    @(214,7,"Terminal Hortofrutícola Agro Chillán","Ñuble",44883,16,100112037,"Cebollín","Sin especificar","Primera",400,600,700,650,"$/paquete 6 unidades","Provincia de Diguillín",108,6,"Hortaliza"),
    @(215,7,"Terminal Hortofrutícola Agro Chillán","Ñuble",44883,16,100112037,"Cebollín","Sin especificar","Segunda",300,500,500,500,"$/paquete 6 unidades","Provincia de Diguillín",83,6,"Hortaliza"),
    @(216,7,"Terminal Hortofrutícola Agro Chillán","Ñuble",44160,16,100112037,"Cebollín","Sin especificar","Primera",43,3500,4000,3709,"$/paquete 36 unidades","Región Metropolitana",103,36,"Hortaliza"),
    @(217,7,"Terminal Hortofrutícola Agro Chillán","Ñuble",44944,16,100112037,"Cebollín","Sin especificar","Primera",400,600,700,650,"$/paquete 6 unidades","Provincia de Diguillín",108,6,"Hortaliza"),
    @(218,7,"Terminal Hortofrutícola Agro Chillán","Ñuble",44944,16,100112037,"Cebollín","Sin especificar","Segunda",300,500,500,500,"$/paquete 6 unidades","Provincia de Diguillín",83,6,"Hortaliza"),
    @(219,7,"Terminal Hortofrutícola Agro Chillán","Ñuble",45016,16,100112037,"Cebollín","Sin especificar","Primera",60,7000,7000,7000,"$/paquete 36 unidades","Provincia de Diguillín",194,36,"Hortaliza"),
    @(220,7,"Terminal Hortofrutícola Agro Chillán","Ñuble",45166,16,100112037,"Cebollín","Sin especificar","Primera",250,6000,6000,6000,"$/paquete 36 unidades","Provincia de Diguillín",167,36,"Hortaliza"),
    @(221,7,"Terminal Hortofrutícola Agro Chillán","Ñuble",45071,16,100112037,"Cebollín","Sin especificar","Primera",170,6000,7000,6529,"$/paquete 36 unidades","Provincia de Diguillín",181,36,"Hortaliza"),
    @(222,7,"Terminal Hortofrutícola Agro Chillán","Ñuble",45223,16,100112037,"Cebollín","Sin especificar","Primera",150,5000,5000,5000,"$/paquete 36 unidades","Provincia de Diguillín",139,36,"Hortaliza"),
    @(223,7,"Terminal Hortofrutícola Agro Chillán","Ñuble",44225,16,100112037,"Cebollín","Sin especificar","Primera",80,3400,3700,3550,"$/paquete 2 kilos","Provincia de Diguillín",1775,2,"Hortaliza"),
    @(224,7,"Terminal Hortofrutícola Agro Chillán","Ñuble",45091,16,100112037,"Cebollín","Sin especificar","Primera",60,5500,5500,5500,"$/paquete 36 unidades","Provincia de Diguillín",153,36,"Hortaliza"),
    @(225,7,"Terminal Hortofrutícola Agro Chillán","Ñuble",45091,16,100112037,"Cebollín","Sin especificar","Segunda",60,4500,4500,4500,"$/paquete 36 unidades","Provincia de Diguillín",125,36,"Hortaliza"),
    @(226,7,"Terminal Hortofrutícola Agro Chillán","Ñuble",45035,16,100112037,"Cebollín","Sin especificar","Primera",150,7000,7000,7000,"$/paquete 36 unidades","Provincia de Diguillín",194,36,"Hortaliza"),
    @(227,7,"Terminal Hortofrutícola Agro Chillán","Ñuble",45121,16,100112037,"Cebollín","Sin especificar","Primera",180,7000,7000,7000,"$/paquete 36 unidades","Provincia de Diguillín",194,36,"Hortaliza"),
    @(228,7,"Terminal Hortofrutícola Agro Chillán","Ñuble",45097,16,100112037,"Cebollín","Sin especificar","Primera",200,6000,7000,6500,"$/paquete 36 unidades","Provincia de Diguillín",181,36,"Hortaliza"),
    @(229,7,"Terminal Hortofrutícola Agro Chillán","Ñuble",45114,16,100112037,"Cebollín","Sin especificar","Primera",80,6000,6000,6000,"$/paquete 36 unidades","Provincia de Diguillín",167,36,"Hortaliza"),
    @(230,7,"Terminal Hortofrutícola Agro Chillán","Ñuble",45239,16,100112037,"Cebollín","Sin especificar","Segunda",60,4500,4500,4500,"$/paquete 36 unidades","Provincia de Diguillín",125,36,"Hortaliza"),
    @(231,7,"Terminal Hortofrutícola Agro Chillán","Ñuble",45173,16,100112037,"Cebollín","Sin especificar","Primera",200,6000,6000,6000,"$/paquete 36 unidades","Provincia de Diguillín",167,36,"Hortaliza"),
    @(232,7,"Terminal Hortofrutícola Agro Chillán","Ñuble",45173,16,100112037,"Cebollín","Sin especificar","Segunda",200,5000,5000,5000,"$/paquete 36 unidades","Provincia de Diguillín",139,36,"Hortaliza")
)

foreach ($r in $rows) {
    $rowNum = $r[0]
    for ($col = 1; $col -le 18; $col++) {
        $ws.Cells.Item($rowNum, $col).Value = $r[$col]
    }
    # Column D (Fecha) must keep the sheet's date number format - this
    # matters for the two newly-appended rows (231-232), which otherwise
    # default to the plain "General" format.
    $ws.Cells.Item($rowNum, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
}
